$wb = $excel.ActiveWorkbook

$wsCases = $wb.Worksheets.Item("testcases")
$wsSteps = $wb.Worksheets.Item("teststeps")

# --- testcases sheet: toggle which test cases run ---
# "role creation" now skipped (its reset case below replaces it)
$wsCases.Activate()
$wsCases.Range("B4").Value = "no"
# "employee creation" now runs
$wsCases.Range("B5").Value = "yes"
$wsCases.Range("B4").Select()

# --- teststeps sheet ---
$wsSteps.Activate()

# complete the "click on reset button" step (row 27) with its
# locator type, locator value and keyword
$wsSteps.Range("C27").Value = "id"
$wsSteps.Range("D27").Value = "Btn_Reset"
$wsSteps.Range("E27").Value = "click"

# update the branch used for employee creation
$wsSteps.Range("F20").Value = "Hyderabad"

# scroll the teststeps view down and land the selection on F20
$excel.ActiveWindow.ScrollRow = 12
$excel.ActiveWindow.ScrollColumn = 1
$wsSteps.Range("F20").Select()
